$wb = $excel.ActiveWorkbook

# 1. Rename the second sheet ("Include from Krebsstadium Cod" -> "Include #0")
$wsInclude = $wb.Worksheets.Item("Include from Krebsstadium Cod")
$wsInclude.Name = "Include #0"

# 2. Update the Metadata sheet
$wsMeta = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B)
$wsMeta.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# Insert a new row before the current "Description" row (row 11), shifting
# Description/Purpose/Copyright/Immutable rows down by one.
# -4121 = xlShiftDown
$wsMeta.Range("A11:B11").Insert(-4121)

# Match the formatting of the surrounding data rows (the insert otherwise
# leaves the new row with default/no formatting).
# -4122 = xlPasteFormats
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)

# Fill the new row with the "Jurisdiction" property (blank value).
$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
